# Scheduled-runner update: refresh cached market-price-derived columns
# (currentAveragePrice[/NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# for a handful of leve rows across the per-job profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1429320.2
$ws.Range("I43").Value = 887
$ws.Range("J43").Value = 3333898
$ws.Range("K43").Value = 887
$ws.Range("L43").Value = 3333898
$ws.Range("M43").Value = -818
$ws.Range("N43").Value = -3334036

$ws.Range("H76").Value = 183339170
$ws.Range("I76").Value = 220005400
$ws.Range("J76").Value = 8000
$ws.Range("K76").Value = 220005400
$ws.Range("L76").Value = 8000
$ws.Range("M76").Value = -220005085
$ws.Range("N76").Value = -8630

$ws.Range("H79").Value = 183339170
$ws.Range("I79").Value = 220005400
$ws.Range("J79").Value = 8000
$ws.Range("K79").Value = 220005400
$ws.Range("L79").Value = 8000
$ws.Range("M79").Value = -220004308
$ws.Range("N79").Value = -10184

$ws.Range("H80").Value = 457915.12
$ws.Range("I80").Value = 443.72726
$ws.Range("K80").Value = 1331.18178
$ws.Range("M80").Value = -333.1817799999999

$ws.Range("H83").Value = 457915.12
$ws.Range("I83").Value = 443.72726
$ws.Range("K83").Value = 3993.54534
$ws.Range("M83").Value = 998.4546599999999

$ws.Range("H107").Value = 769347.4399999999
$ws.Range("I107").Value = 833434.75
$ws.Range("K107").Value = 833434.75
$ws.Range("M107").Value = -831514.75

$ws.Range("H132").Value = 628330.0600000001
$ws.Range("I132").Value = 4555.5557
$ws.Range("J132").Value = 1430325.9
$ws.Range("K132").Value = 13666.6671
$ws.Range("L132").Value = 4290977.699999999
$ws.Range("M132").Value = -11136.6671
$ws.Range("N132").Value = -4296037.699999999

$ws.Range("H137").Value = 902.65
$ws.Range("I137").Value = 816.4375
$ws.Range("J137").Value = 1247.5
$ws.Range("K137").Value = 2449.3125
$ws.Range("L137").Value = 3742.5
$ws.Range("M137").Value = 100.6875
$ws.Range("N137").Value = -8842.5

$ws.Range("H138").Value = 3138.72
$ws.Range("I138").Value = 786.5238000000001
$ws.Range("J138").Value = 3763.9873
$ws.Range("K138").Value = 2359.5714
$ws.Range("L138").Value = 11291.9619
$ws.Range("M138").Value = 2780.4286
$ws.Range("N138").Value = -21571.9619

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1854.6923
$ws.Range("I45").Value = 1685.1666
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1685.1666
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1308.1666
$ws.Range("N45").Value = -2754

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 142943890
$ws.Range("I20").Value = 150975
$ws.Range("J20").Value = 333334430
$ws.Range("K20").Value = 150975
$ws.Range("L20").Value = 333334430
$ws.Range("M20").Value = -150728
$ws.Range("N20").Value = -333334924

$ws.Range("H134").Value = 38655.816
$ws.Range("I134").Value = 1588.28
$ws.Range("J134").Value = 502000
$ws.Range("K134").Value = 4764.84
$ws.Range("L134").Value = 1506000
$ws.Range("M134").Value = -2229.84
$ws.Range("N134").Value = -1511070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1527.4
$ws.Range("I16").Value = 1516.2307
$ws.Range("J16").Value = 1600
$ws.Range("K16").Value = 1516.2307
$ws.Range("L16").Value = 1600
$ws.Range("M16").Value = -1229.2307
$ws.Range("N16").Value = -2174

$ws.Range("H31").Value = 1655.8125
$ws.Range("I31").Value = 1414.3
$ws.Range("J31").Value = 2058.3333
$ws.Range("K31").Value = 1414.3
$ws.Range("L31").Value = 2058.3333
$ws.Range("M31").Value = -1119.3
$ws.Range("N31").Value = -2648.3333

$ws.Range("H34").Value = 1655.8125
$ws.Range("I34").Value = 1414.3
$ws.Range("J34").Value = 2058.3333
$ws.Range("K34").Value = 1414.3
$ws.Range("L34").Value = 2058.3333
$ws.Range("M34").Value = -1212.3
$ws.Range("N34").Value = -2462.3333

$ws.Range("H99").Value = 2770.75
$ws.Range("I99").Value = 2487.4375
$ws.Range("J99").Value = 3337.375
$ws.Range("K99").Value = 2487.4375
$ws.Range("L99").Value = 3337.375
$ws.Range("M99").Value = -989.4375
$ws.Range("N99").Value = -6333.375

$ws.Range("H113").Value = 1527.4
$ws.Range("I113").Value = 1516.2307
$ws.Range("J113").Value = 1600
$ws.Range("K113").Value = 1516.2307
$ws.Range("L113").Value = 1600
$ws.Range("M113").Value = 653.7692999999999
$ws.Range("N113").Value = -5940

$ws.Range("H126").Value = 2770.75
$ws.Range("I126").Value = 2487.4375
$ws.Range("J126").Value = 3337.375
$ws.Range("K126").Value = 7462.3125
$ws.Range("L126").Value = 10012.125
$ws.Range("M126").Value = -4992.3125
$ws.Range("N126").Value = -14952.125

$ws.Range("H134").Value = 2279.0356
$ws.Range("I134").Value = 1500.5416
$ws.Range("J134").Value = 6950
$ws.Range("K134").Value = 4501.6248
$ws.Range("L134").Value = 20850
$ws.Range("M134").Value = -1966.6248
$ws.Range("N134").Value = -25920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 600
$ws.Range("J32").Value = 600
$ws.Range("L32").Value = 1800
$ws.Range("N32").Value = -2366

$ws.Range("H46").Value = 5000
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 15000
$ws.Range("N46").Value = -15182

$ws.Range("H68").Value = 1025.5
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1025.5
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4581.2666
$ws.Range("I70").Value = 4093.7693
$ws.Range("J70").Value = 7750
$ws.Range("K70").Value = 4093.7693
$ws.Range("L70").Value = 7750
$ws.Range("M70").Value = -3823.7693
$ws.Range("N70").Value = -8290

$ws.Range("H73").Value = 4581.2666
$ws.Range("I73").Value = 4093.7693
$ws.Range("J73").Value = 7750
$ws.Range("K73").Value = 4093.7693
$ws.Range("L73").Value = 7750
$ws.Range("M73").Value = -3157.7693
$ws.Range("N73").Value = -9622

$ws.Range("H113").Value = 1875
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 670

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 491.95456
$ws.Range("I46").Value = 483.1
$ws.Range("J46").Value = 499.33334
$ws.Range("K46").Value = 483.1
$ws.Range("L46").Value = 499.33334
$ws.Range("M46").Value = -295.1
$ws.Range("N46").Value = -875.33334

$ws.Range("H125").Value = 34163.332
$ws.Range("J125").Value = 34163.332
$ws.Range("L125").Value = 34163.332
$ws.Range("N125").Value = -44003.332

$ws.Range("H127").Value = 44612.5
$ws.Range("J127").Value = 44612.5
$ws.Range("L127").Value = 44612.5
$ws.Range("N127").Value = -54532.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 229.48276
$ws.Range("I113").Value = 221.75
$ws.Range("J113").Value = 266.6
$ws.Range("K113").Value = 665.25
$ws.Range("L113").Value = 799.8000000000001
$ws.Range("M113").Value = 1504.75
$ws.Range("N113").Value = -5139.8

$ws.Range("H131").Value = 47512
$ws.Range("J131").Value = 47512
$ws.Range("L131").Value = 47512
$ws.Range("N131").Value = -57592
